$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "1. EXECUTIVE SUMMARY" Heading1 paragraph (the report-body one, not
#    the TOC entry nor the appendix heading) -- apply numbered-list
#    formatting (numId 10 / ilvl 0), drop the manual "1. " prefix, and
#    insert a blank paragraph right after it (before the KPI table).
# ---------------------------------------------------------------------
$execHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "1. EXECUTIVE SUMMARY`r" -and $p.Style.NameLocal -eq "Heading 1") {
        $execHeading = $p
        break
    }
}

$execHeading.Range.ListFormat.ApplyNumberDefault()
$lt = $execHeading.Range.ListFormat.ListTemplate
$levelNumberStyles = @(0, 4, 2, 0, 4, 2, 0, 4, 2)
for ($i = 1; $i -le 9; $i++) {
    $lvl = $lt.ListLevels.Item($i)
    $lvl.NumberStyle = $levelNumberStyles[$i - 1]
    $lvl.NumberFormat = "%$i."
}

$found = $execHeading.Range.Find.Execute("1. EXECUTIVE SUMMARY", $true, $false, $false, $false, $false, $true, 1, $false, "EXECUTIVE SUMMARY", 2)

$execHeading.Range.InsertParagraphAfter()
$blankAfterExec = $execHeading.Next()
$blankAfterExec.Range.Delete()
$blankAfterExec.Style = $d.Styles.Item("Normal")

# ---------------------------------------------------------------------
# 2) "6. BANK-SPECIFIC ANALYSIS" Heading1 paragraph -- just insert a
#    blank paragraph right after it (before the bank comparison table).
# ---------------------------------------------------------------------
$bankHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "6. BANK-SPECIFIC ANALYSIS`r" -and $p.Style.NameLocal -eq "Heading 1") {
        $bankHeading = $p
        break
    }
}

$bankHeading.Range.InsertParagraphAfter()
$blankAfterBank = $bankHeading.Next()
$blankAfterBank.Range.Delete()
$blankAfterBank.Style = $d.Styles.Item("Normal")

# ---------------------------------------------------------------------
# 3) Merge the leading-space run with the "Opportunities: ..." run in
#    the 6.1 bullet list into a single run.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(" Opportunities: Error handling improvements, UX simplification", $true, $false, $false, $false, $false, $true, 1, $false, " Opportunities: Error handling improvements, UX simplification", 2)
